$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the data between row 3 and row 4 for the columns that differ
# (D: Fecha, J: Volumen, K: Precio minimo, L: Precio maximo, M: Precio promedio ponderado, P: Precio $/Kg)

$row3 = @{
    D = $ws.Range("D3").Value2
    J = $ws.Range("J3").Value2
    K = $ws.Range("K3").Value2
    L = $ws.Range("L3").Value2
    M = $ws.Range("M3").Value2
    P = $ws.Range("P3").Value2
}

$row4 = @{
    D = $ws.Range("D4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}

$ws.Range("D3").Value = $row4.D
$ws.Range("J3").Value = $row4.J
$ws.Range("K3").Value = $row4.K
$ws.Range("L3").Value = $row4.L
$ws.Range("M3").Value = $row4.M
$ws.Range("P3").Value = $row4.P

$ws.Range("D4").Value = $row3.D
$ws.Range("J4").Value = $row3.J
$ws.Range("K4").Value = $row3.K
$ws.Range("L4").Value = $row3.L
$ws.Range("M4").Value = $row3.M
$ws.Range("P4").Value = $row3.P
